# Update the deal_value column (D) with the new probability-based results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 937.5377473015355
    3  = 937.5377473015355
    4  = 996.2938725423382
    5  = 1116.460497103601
    6  = 1240.349233061726
    7  = 1368.198052860202
    8  = 1500.268506819392
    9  = 1636.848944567762
    10 = 1778.258306167385
    11 = 1924.850608277711
    12 = 1116.460497103601
    13 = 1116.460497103601
    14 = 1116.460497103601
    15 = 1116.460497103601
    16 = 1240.349233061726
    17 = 1368.198052860202
    18 = 1500.268506819392
    19 = 1636.848944567762
    20 = 1778.258306167385
    21 = 1924.850608277711
    22 = 1116.460497103601
    23 = 1116.460497103601
    24 = 1116.460497103601
    25 = 1116.460497103601
    26 = 1240.349233061726
    27 = 1368.198052860202
    28 = 1500.268506819392
    29 = 1636.848944567762
    30 = 1778.258306167385
    31 = 1924.850608277711
}

foreach ($row in $newValues.Keys) {
    $ws.Range("D$row").Value = $newValues[$row]
}
